$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (was "max"/"prediction" derived numeric), per row 2..8
$values = @(1610.094102676912, 1864.203183037275, 1270.202934235216, 1728.49172040795, 1639.053738423656, 1588.61140049857, 1637.953737106925)

# Delete column C ("max") entirely, shifting D (prediction) and E (rejection-f) left
$ws.Columns.Item(3).Delete()

# Update column B values (now numeric prediction-distance values instead of 1s)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Column C (previously D, "prediction") now holds inline strings instead of numbers - already string typed in original, keep as is.
